$d = $word.ActiveDocument

# The "Dauer" (duration) paragraph currently reads "15min". It needs to
# become "25min", with the existing "_GoBack" bookmark repositioned right
# after the "2" -- i.e. the text is split into two runs ("2" and "5min")
# with the bookmark sitting between them, just like Word leaves it after a
# manual edit at the caret.
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "15min*") {
        $target = $p
        break
    }
}

if ($target -ne $null) {
    $target.Range.Text = "25min"

    $start = $target.Range.Start
    $bmRange = $d.Range($start + 1, $start + 1)
    $d.Bookmarks.Add("_GoBack", $bmRange)
}

$d.Save()
